# "day ten first commit"
# Adds a second sheet ("Concepts to remeber"), a defined name "Leetcode",
# five new rows of question-tracker data (rows 5-9) on "Questions Tracker",
# plus matching hyperlinks/formatting, and populates the new sheet with
# "concepts" notes.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------------
# 1. Workbook-level: defined name + second worksheet
# ---------------------------------------------------------------------------
$wb.Names.Add("Leetcode", '''Questions Tracker''!$C$5')

$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Concepts to remeber"

# ===========================================================================
# 2. "Questions Tracker" sheet - header tweak + 5 new rows
# ===========================================================================

# --- header cell C1 gets center/middle/wrap like the rest of the header ---
$ws1.Range("C1").HorizontalAlignment = -4108   # xlCenter
$ws1.Range("C1").VerticalAlignment = -4108     # xlCenter
$ws1.Range("C1").WrapText = $true

# --- existing C2:C4 (empty "Link" cells) get center/middle/wrap border style ---
$ws1.Range("C2:C4").HorizontalAlignment = -4108
$ws1.Range("C2:C4").VerticalAlignment = -4108
$ws1.Range("C2:C4").WrapText = $true

# --- column widths (A,B,C,D,E,F tightened up a touch; H widened) ---
$ws1.Columns.Item(1).ColumnWidth = 19.08984375
$ws1.Columns.Item(2).ColumnWidth = 35.54296875
$ws1.Columns.Item(3).ColumnWidth = 26.26953125
$ws1.Columns.Item(5).ColumnWidth = 16.08984375
$ws1.Columns.Item(6).ColumnWidth = 12
$ws1.Columns.Item(8).ColumnWidth = 40.81640625

# ---------------------------------------------------------------------------
# Row 5 - While loop / 1281. Subtract the Product and Sum of Digits ...
# ---------------------------------------------------------------------------
$ws1.Range("A2").Copy()
$ws1.Range("A5").PasteSpecial(-4122)
$ws1.Range("A5").Value = "While loop"

$ws1.Range("B2").Copy()
$ws1.Range("B5").PasteSpecial(-4122)
$ws1.Range("B5").WrapText = $true
$ws1.Range("B5").Value = "1281. Subtract the Product and Sum of Digits of an Integer"

$ws1.Range("C2").Copy()
$ws1.Range("C5").PasteSpecial(-4122)
$ws1.Range("C5").HorizontalAlignment = -4108
$ws1.Range("C5").VerticalAlignment = -4108
$ws1.Range("C5").WrapText = $true
$ws1.Range("C5").Value = "LeetCode"
$ws1.Hyperlinks.Add($ws1.Range("C5"), "https://leetcode.com/problems/subtract-the-product-and-sum-of-digits-of-an-integer/") | Out-Null

$ws1.Range("D2").Copy()
$ws1.Range("D5").PasteSpecial(-4122)
$ws1.Range("D5").Value = 1

$ws1.Range("E2").Copy()
$ws1.Range("E5").PasteSpecial(-4122)
$ws1.Range("E5").Value = 44845

$ws1.Range("F2").Copy()
$ws1.Range("F5").PasteSpecial(-4122)
$ws1.Range("F5").Value = 0

$ws1.Range("G2").Copy()
$ws1.Range("G5").PasteSpecial(-4122)
$ws1.Range("G5").Value = "No"

$ws1.Range("H2").Copy()
$ws1.Range("H5").PasteSpecial(-4122)
$ws1.Range("H5").Value = "Basic Question"

$ws1.Range("A5:H5").RowHeight = 29

# ---------------------------------------------------------------------------
# Row 6 - Bitwise Operation / 191. Number of 1 Bits
# ---------------------------------------------------------------------------
$ws1.Range("A5").Copy()
$ws1.Range("A6").PasteSpecial(-4122)
$ws1.Range("A6").Value = "Bitwise Operation"

$ws1.Range("B2").Copy()
$ws1.Range("B6").PasteSpecial(-4122)
$ws1.Range("B6").Value = "191. Number of 1 Bits"

$ws1.Range("C5").Copy()
$ws1.Range("C6").PasteSpecial(-4122)
$ws1.Range("C6").Value = "LeetCode"
$ws1.Hyperlinks.Add($ws1.Range("C6"), "https://leetcode.com/problems/number-of-1-bits/") | Out-Null

$ws1.Range("D5").Copy()
$ws1.Range("D6").PasteSpecial(-4122)
$ws1.Range("D6").Value = 1

$ws1.Range("E5").Copy()
$ws1.Range("E6").PasteSpecial(-4122)
$ws1.Range("E6").Value = 44845

$ws1.Range("F5").Copy()
$ws1.Range("F6").PasteSpecial(-4122)
$ws1.Range("F6").Value = 1

$ws1.Range("G2").Copy()
$ws1.Range("G6").PasteSpecial(-4122)
$ws1.Range("G6").Value = "Yes"

$ws1.Range("H5").Copy()
$ws1.Range("H6").PasteSpecial(-4122)
$ws1.Range("H6").Value = "Basic Question"

$ws1.Range("A6:H6").RowHeight = 18

# ---------------------------------------------------------------------------
# Row 7 - Bitwise Operation / Decimal to Binary conversion
# ---------------------------------------------------------------------------
$ws1.Range("A6").Copy()
$ws1.Range("A7").PasteSpecial(-4122)
$ws1.Range("A7").Value = "Bitwise Operation"

$ws1.Range("B6").Copy()
$ws1.Range("B7").PasteSpecial(-4122)
$ws1.Range("B7").Value = "Decimal to Binary conversion"

$ws1.Range("C2").Copy()
$ws1.Range("C7").PasteSpecial(-4122)

$ws1.Range("D6").Copy()
$ws1.Range("D7").PasteSpecial(-4122)
$ws1.Range("D7").Value = 1

$ws1.Range("E6").Copy()
$ws1.Range("E7").PasteSpecial(-4122)
$ws1.Range("E7").Value = 44845

$ws1.Range("F6").Copy()
$ws1.Range("F7").PasteSpecial(-4122)
$ws1.Range("F7").Value = 0

$ws1.Range("G6").Copy()
$ws1.Range("G7").PasteSpecial(-4122)
$ws1.Range("G7").Value = "Yes"

$ws1.Range("H6").Copy()
$ws1.Range("H7").PasteSpecial(-4122)
$ws1.Range("H7").Value = "Basic Question"

$ws1.Range("A7:H7").RowHeight = 18

# ---------------------------------------------------------------------------
# Row 8 - Bitwise Operation / Negative number to binary conversion
# ---------------------------------------------------------------------------
$ws1.Range("A7").Copy()
$ws1.Range("A8").PasteSpecial(-4122)
$ws1.Range("A8").Value = "Bitwise Operation"

$ws1.Range("B7").Copy()
$ws1.Range("B8").PasteSpecial(-4122)
$ws1.Range("B8").Value = "Negative number to binary conversion"

$ws1.Range("C7").Copy()
$ws1.Range("C8").PasteSpecial(-4122)

$ws1.Range("D7").Copy()
$ws1.Range("D8").PasteSpecial(-4122)
$ws1.Range("D8").Value = 1

$ws1.Range("E7").Copy()
$ws1.Range("E8").PasteSpecial(-4122)
$ws1.Range("E8").Value = 44876

$ws1.Range("F7").Copy()
$ws1.Range("F8").PasteSpecial(-4122)
$ws1.Range("F8").Value = 2

$ws1.Range("G7").Copy()
$ws1.Range("G8").PasteSpecial(-4122)
$ws1.Range("G8").Value = "Yes"

$ws1.Range("H7").Copy()
$ws1.Range("H8").PasteSpecial(-4122)
$ws1.Range("H8").Value = "Basic Question"

# Row 8 itself keeps default row height, but the font/fill need to match a
# "no alignment" variant - nudge via explicit Borders reapply + Fill
$ws1.Range("A8,B8,C8,D8,E8,F8,G8,H8").WrapText = $false

# ---------------------------------------------------------------------------
# Row 9 - Program for replacing one digit with other / geeksforgeeks
# ---------------------------------------------------------------------------
$ws1.Range("A8").Copy()
$ws1.Range("A9").PasteSpecial(-4122)
$ws1.Range("A9").ClearContents()

$ws1.Range("B8").Copy()
$ws1.Range("B9").PasteSpecial(-4122)
$ws1.Range("B9").VerticalAlignment = -4108
$ws1.Range("B9").Value = "Program for replacing one digit with other"

$ws1.Range("C5").Copy()
$ws1.Range("C9").PasteSpecial(-4122)
$ws1.Range("C9").Value = "geeksforgeeks"
$ws1.Hyperlinks.Add($ws1.Range("C9"), "https://www.geeksforgeeks.org/") | Out-Null

$ws1.Range("D8").Copy()
$ws1.Range("D9").PasteSpecial(-4122)
$ws1.Range("D9").Value = 1

$ws1.Range("E8").Copy()
$ws1.Range("E9").PasteSpecial(-4122)
$ws1.Range("E9").Value = 44876

$ws1.Range("F8").Copy()
$ws1.Range("F9").PasteSpecial(-4122)
$ws1.Range("F9").Value = 0

$ws1.Range("G8").Copy()
$ws1.Range("G9").PasteSpecial(-4122)
$ws1.Range("G9").Value = "Yes"

$ws1.Range("H8").Copy()
$ws1.Range("H9").PasteSpecial(-4122)
$ws1.Range("H9").WrapText = $true
$ws1.Range("H9").Value = "Having problem in 2's compliment as how to replace just a single digit in integer so googled it and found this solution"

$ws1.Range("A9:H9").RowHeight = 43.5

$ws1.Range("F16").Select()

# ===========================================================================
# 3. "Concepts to remeber" sheet
# ===========================================================================
$ws2.Columns.Item(1).ColumnWidth = 31.54296875
$ws2.Columns.Item(2).ColumnWidth = 58.1796875
$ws2.Columns.Item(3).ColumnWidth = 37.81640625

# ---- Row 1: headers ----
$ws2.Range("A1").Value = "Concept"
$ws2.Range("B1").Value = "Details"
$ws2.Range("C1").Value = "Lecture covered"
$ws2.Range("A1:C1").Font.Bold = $true
$ws2.Range("A1:C1").HorizontalAlignment = -4131
$ws2.Range("A1:C1").VerticalAlignment = -4160
$ws2.Range("B1").WrapText = $true

$ws2.Range("A1:C1").Borders.LineStyle = 1
$ws2.Range("A1:C1").Borders.Weight = 2

# ---- Row 2 ----
$ws2.Range("A2").Value = "number & 1"
$ws2.Range("B2").Value = "When we do AND of any number with 1 and if it returns 0 then it is an odd number else it is an even number.`nNumber&1 - basically does the AND of last bit of that number with 1."
$ws2.Range("C2").Value = "Lecture 6: Binary & Decimal Number System"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://www.geeksforgeeks.org/") | Out-Null
$ws2.Range("A2:C2").RowHeight = 59

# ---- Row 3 ----
$ws2.Range("A3").Value = "answer = (digit * answer^i) + answer"
$ws2.Range("B3").Value = "Formula to reverse an integer (ex: 123 -> 321)"
$ws2.Range("C3").Value = "Lecture 6: Binary & Decimal Number System"
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://www.geeksforgeeks.org/") | Out-Null
$ws2.Range("A3:C3").RowHeight = 30

# ---- Row 4 ----
$ws2.Range("A4").Value = "answer = (10* answer) + digit"
$ws2.Range("B4").Value = "Formula to print integer in same flow (ex: 123 -> 123)"
$ws2.Range("C4").Value = "Lecture 6: Binary & Decimal Number System"
$ws2.Hyperlinks.Add($ws2.Range("C4"), "https://www.geeksforgeeks.org/") | Out-Null
$ws2.Range("A4:C4").RowHeight = 27.5

# ---- Row 5 ----
$ws2.Range("A5").Value = "Integer Range -> (-2^31) to (2^31-1)"
$ws2.Range("B5").Value = "If you get out of this range it will return garbage value"
$ws2.Range("C5").Value = ""

$ws2.Range("A2:A5,B2:B5").VerticalAlignment = -4160
$ws2.Range("A2:A5,B2:B5").HorizontalAlignment = -4131
$ws2.Range("B2:B5").WrapText = $true

$ws2.Range("A1:C5").Borders.LineStyle = 1
$ws2.Range("A1:C5").Borders.Weight = 2

$ws2.Range("B8").Select()
